$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C20 (EXCLUDE.PARRUN row): fill in example value "1-2"
$ws.Range("C20").Value = "1-2"

# Row 22: new "Motion" category entry for MTN.OVERWRITE
$ws.Range("A22").Value = "Motion"
$ws.Range("B22").Value = "Create motion plots even if already run (overwrites prior)"
$ws.Range("C22").Value = $true
$ws.Range("D22").Value = "TRUE or FALSE. Set true to always generate new motion plots."
$ws.Range("E22").Value = "MTN.OVERWRITE"

# Move the view so row 6 is at the top and C22 is the active selection
$ws.Range("C22").Select()
$excel.ActiveWindow.ScrollRow = 6

$wb.Save()
